$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data: Date serial, y_0 (year), y_0_forecast, y_1 (year), y_1_forecast
$data = @(
    @(39400, 2007, 1.144978573787081,   2008, 1.854677717601594),
    @(39765, 2008, 1.381024225294869,   2009, 1.929644353290927),
    @(40130, 2009, -0.8792832172735965, 2010, -0.10374497415091),
    @(40494, 2010, 0.9337833426867448,  2011, 1.899713704008654),
    @(40862, 2011, 2.791140000794279,   2012, 1.722110645261998),
    @(41228, 2012, 0.4451370000809973,  2013, 0.6480763427741953),
    @(41592, 2013, 0.2545814083968478,  2014, 1.223618887196531),
    @(41957, 2014, 1.297015177357297,   2015, 0.309136676902555),
    @(42321, 2015, 1.365576377841027,   2016, 1.582979977679533),
    @(42689, 2016, 2.204449574611278,   2017, 1.401113624217043),
    @(43053, 2017, 2.18621550610123,    2018, 2.357704431248409),
    @(43418, 2018, 1.911050033324102,   2019, 2.362613045431528),
    @(43783, 2019, 1.457852003181337,   2020, 0.7756897792100315),
    @(44159, 2020, -3.258619210312885,  2021, -2.035749133083187),
    @(44525, 2021, 0.4255262881966981,  2022, 2.107044281672787),
    @(44890, 2022, 3.293290997728171,   2023, 1.50842747477733),
    @(45254, 2023, -0.2814561130375925, 2024, 0.06837453075889677),
    @(45618, 2024, -0.6470065423293869, 2025, 1.000784331770643)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}

# Row 19 is new; make sure column A keeps the same style as the rest of column A (date style)
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
